$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new blank row at the top of the data (row 2), shifting rows 2-15 down to 3-16 ---
$ws.Rows("2:2").Insert()

# --- 2. Copy cell formatting (number format / alignment / font / border) from row 3 onto the
#        freshly inserted row 2 so the new row matches the look of the rest of the table
#        (row-insert otherwise invents a bold "header-like" style for the new row). ---
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- 3. Populate the new row 2 with the latest price entry ---
$ws.Range("A2").Value = 15
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 294.05

# Force the circular-date column to be stored as literal text (not auto-parsed into a date
# serial number), then restore its format to match the sibling cells in column E.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "30-10-2025"
$ws.Range("E3").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("F2").Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf"

# --- 4. Rebuild all hyperlinks in column F (rows 2-16). The row-insert operation shifts cell
#        text/values but leaves old hyperlink anchors pointing at their original rows, so the
#        safest approach is to drop every hyperlink in the column and re-create them fresh,
#        matching each row's (now correct) URL text. ---
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 6).Hyperlinks.Delete()
}

$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(15, 6), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Cells.Item(16, 6), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf")

# --- 5. Hyperlinks.Add() applies Excel's built-in blue/underlined "Hyperlink" style to the
#        target cells, but the source workbook keeps column F in the same plain centered style
#        as the rest of the table. Re-apply the plain style (copied from column A, which carries
#        the same base alignment/font/border) onto column F without touching the cell text. ---
$ws.Range("A2:A16").Copy()
$ws.Range("F2:F16").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
